$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.728.75"
$ws.Range("E2").Value = "  -2.10%  "
$ws.Range("D3").Value = "3.087.71"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.74"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.30"
$ws.Range("E6").Value = "  -2.79%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "3.087.99"
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.14"
$ws.Range("E10").Value = "  -3.30%  "
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.391"
$ws.Range("E12").Value = "  +2.19%  "
$ws.Range("D13").Value = "3.622.73"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.45"
$ws.Range("E15").Value = "  -6.11%  "
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("D17").Value = "57.812.10"
$ws.Range("E17").Value = "  -1.94%  "
$ws.Range("D18").Value = "3.101.60"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.65"
$ws.Range("E20").Value = "  -3.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.93"
$ws.Range("E21").Value = "  -3.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "341.57"
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.510"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.30"
$ws.Range("E25").Value = "  +2.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.169"
$ws.Range("E26").Value = "  -1.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").Value = "0.0₃0911"
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.36"
$ws.Range("E30").Value = "  -5.46%  "
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.86"
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.90"
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.17"
$ws.Range("E34").Value = "  -3.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.67"
$ws.Range("E35").Value = "  +1.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.60"
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.12"
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.92"
$ws.Range("E38").Value = "  -4.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.23"
$ws.Range("E39").Value = "  -5.83%  "
$ws.Range("E40").Value = "  -3.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.99"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("E42").Value = "  +5.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.683"
$ws.Range("E43").Value = "  +2.38%  "
$ws.Range("D44").Value = "3.130.63"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.89"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("E47").Value = "  +1.58%  "
$ws.Range("D48").Value = "2.276.98"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.986"
$ws.Range("E49").Value = "  +2.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.08"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.36"
$ws.Range("E51").Value = "  -3.48%  "
